$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-07-19 Saturday" "2025-07-20 Sunday"
Replace-Text "399×2=" "505×3="
Replace-Text "888×3=" "552×8="
Replace-Text "667×6=" "286×2="
Replace-Text "957×5=" "757×9="
Replace-Text "775×6=" "795×3="
Replace-Text "555×2=" "827×8="
Replace-Text "444×9=" "255×2="
Replace-Text "162×4=" "399×8="
Replace-Text "972×2=" "223×6="
Replace-Text "469×9=" "621×2="
Replace-Text "312×9=" "209×4="
Replace-Text "870×4=" "527×2="
Replace-Text "816×2=" "192×5="
Replace-Text "646×5=" "699×8="
Replace-Text "156×7=" "510×8="
Replace-Text "264×9=" "764×5="
Replace-Text "786×3=" "414×3="
Replace-Text "592×8=" "528×7="
Replace-Text "781×7=" "156×7="
Replace-Text "542×9=" "623×4="
Replace-Text "809×8=" "197×8="
Replace-Text "905×2=" "367×8="
Replace-Text "268×5=" "453×9="
Replace-Text "761×8=" "120×6="
Replace-Text "565×7=" "257×6="
